# Apply the "new version with timestamp" update to the DaySale report.
# Updates the balance/sale-price/transaction-count figures for a few
# products (stock recount + extra sales recorded), refreshes the running
# total, and bumps the generated-at timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DEXAMETHASONE-MUP 8MG/2ML 5 AMP (row 20)
$ws.Range("H20").Value = "3:0"
$ws.Range("P20").Value = "26.0000"
$ws.Range("Q20").Value = "0:2"

# VOLTAREN 75MG/3ML 3 AMP. (row 32)
$ws.Range("H32").Value = "1:3"
$ws.Range("P32").Value = "33.6600"
$ws.Range("Q32").Value = "0:2"

# سرنجات 3 سم (row 37)
$ws.Range("P37").Value = "8.0000"
$ws.Range("Q37").Value = "4:0"

# سرنجات 5 سم (row 38)
$ws.Range("P38").Value = "9.0000"
$ws.Range("Q38").Value = "3:0"

# Refresh the running total of the "sale price" column
$ws.Range("P42").Value = 1746.385

# Bump the "generated at" timestamp in the footer
$ws.Range("A43").Value = "Thursday, 21 August, 2025 1:10 PM"
